# Update "想去人数" (interested-people count) values in column F across
# the workbook's four sheets, per the commit "Update gh-pages to output
# generated at 456a3b4".

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$sheet1Updates = @{
    4  = 635
    5  = 473
    6  = 260
    7  = 1144
    9  = 175
    10 = 64
    11 = 763
    12 = 407
    15 = 199
    17 = 383
    18 = 6230
    20 = 56
    22 = 7202
    23 = 26
    24 = 27
    25 = 3291
    26 = 415
    27 = 791
    28 = 4475
    29 = 331
    30 = 155
    32 = 1267
    33 = 116
    34 = 37
    36 = 990
    37 = 1287
    38 = 2073
}
foreach ($row in $sheet1Updates.Keys) {
    $ws1.Range("F$row").Value = $sheet1Updates[$row]
}

# Sheet "演出" (Performance)
$ws2 = $wb.Worksheets.Item("演出")
$sheet2Updates = @{
    3 = 34
    5 = 70
}
foreach ($row in $sheet2Updates.Keys) {
    $ws2.Range("F$row").Value = $sheet2Updates[$row]
}

# Sheet "本地生活" (Local life)
$ws3 = $wb.Worksheets.Item("本地生活")
$sheet3Updates = @{
    3 = 1168
    4 = 59
}
foreach ($row in $sheet3Updates.Keys) {
    $ws3.Range("F$row").Value = $sheet3Updates[$row]
}

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$sheet4Updates = @{
    4  = 1168
    5  = 59
    7  = 635
    8  = 473
    9  = 260
    10 = 1144
    12 = 175
    13 = 64
    14 = 763
    15 = 407
    19 = 199
    21 = 383
    22 = 6230
    23 = 6230
    25 = 56
    27 = 7202
    28 = 26
    29 = 27
    30 = 3291
    31 = 415
    32 = 791
    33 = 4475
    34 = 331
    35 = 34
    36 = 155
    38 = 1267
    39 = 116
    40 = 37
    42 = 990
    43 = 1287
    45 = 2073
    47 = 70
}
foreach ($row in $sheet4Updates.Keys) {
    $ws4.Range("F$row").Value = $sheet4Updates[$row]
}
